# Update generated stats (想去人数 / interest counts) on the "展览" and
# "全部类型" sheets, matching the data refresh captured in the commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 406
$ws1.Range("F5").Value = 8468
$ws1.Range("F7").Value = 10525
$ws1.Range("F23").Value = 64
$ws1.Range("F28").Value = 577
$ws1.Range("F30").Value = 1150
$ws1.Range("F43").Value = 631

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 406
$ws4.Range("F9").Value = 8468
$ws4.Range("F11").Value = 10525
$ws4.Range("F26").Value = 577
$ws4.Range("F28").Value = 1150
$ws4.Range("F47").Value = 631
